$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the two new columns' headers: I1 = "I0", J1 = "IF" ---
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the bold / bordered / centered-top formatting already used across row 1
# by copying the format from the existing "IP" header cell (H1) onto the new ones.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# --- Fill in the I0 / IF data values for rows 2-60 ---
# Each tuple is (row, I-value, J-value)
$data = @(@(2,8,8),@(3,8,8),@(4,8,8),@(5,9,10),@(6,7,7),@(7,8,8),@(8,7,7),@(9,6,6),@(10,6,7),@(11,4,5),@(12,7,7),@(13,9,9),@(14,10,10),@(15,7,7),@(16,7,7),@(17,6,6),@(18,7,8),@(19,9,9),@(20,9,9),@(21,9,9),@(22,7,8),@(23,7,7),@(24,8,8),@(25,10,10),@(26,7,7),@(27,9,9),@(28,10,10),@(29,9,9),@(30,8,8),@(31,6,6),@(32,8,9),@(33,7,7),@(34,7,7),@(35,9,9),@(36,9,9),@(37,6,6),@(38,12,12),@(39,6,6),@(40,7,7),@(41,8,8),@(42,8,8),@(43,9,9),@(44,9,9),@(45,6,6),@(46,6,6),@(47,7,8),@(48,8,8),@(49,6,6),@(50,8,8),@(51,6,6),@(52,5,6),@(53,7,7),@(54,7,7),@(55,7,7),@(56,8,8),@(57,7,7),@(58,5,6),@(59,5,5),@(60,4,4))

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 9).Value = $row[1]
    $ws.Cells.Item($r, 10).Value = $row[2]
}
